$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41; this shifts rows 41..157 down to 42..158
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new data point
$ws.Cells.Item(41, 1).Value = 5
$ws.Cells.Item(41, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(41, 3).Value = "Maule"
$ws.Cells.Item(41, 4).Value = 44497
$ws.Cells.Item(41, 5).Value = 7
$ws.Cells.Item(41, 6).Value = 100112008
$ws.Cells.Item(41, 7).Value = "Coliflor"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 5000
$ws.Cells.Item(41, 11).Value = 600
$ws.Cells.Item(41, 12).Value = 600
$ws.Cells.Item(41, 13).Value = 600
$ws.Cells.Item(41, 14).Value = "$/unidad"
$ws.Cells.Item(41, 15).Value = "Región del Maule"
$ws.Cells.Item(41, 16).Value = 600
$ws.Cells.Item(41, 17).Value = 1
$ws.Cells.Item(41, 18).Value = "Hortaliza"
